$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# hunk 0: ALC row 17
$ws1.Range("H17").Value = 870709.75
$ws1.Range("J17").Value = 1044756.75
$ws1.Range("L17").Value = 3134270.25
$ws1.Range("N17").Value = -3134606.25

# hunk 1: ALC row 40
$ws1.Range("H40").Value = 9578.261
$ws1.Range("I40").Value = 6020.1816
$ws1.Range("J40").Value = 12839.833
$ws1.Range("K40").Value = 6020.1816
$ws1.Range("L40").Value = 12839.833
$ws1.Range("M40").Value = -5845.1816
$ws1.Range("N40").Value = -13189.833

# hunk 2: ALC row 64
$ws1.Range("H64").Value = 9266.767
$ws1.Range("J64").Value = 9266.767
$ws1.Range("L64").Value = 9266.767
$ws1.Range("N64").Value = -9762.767

# hunk 3: ALC row 67
$ws1.Range("H67").Value = 9266.767
$ws1.Range("J67").Value = 9266.767
$ws1.Range("L67").Value = 9266.767
$ws1.Range("N67").Value = -10982.767

# hunk 4: ALC row 137
$ws1.Range("H137").Value = 2666.2654
$ws1.Range("I137").Value = 1455.0952
$ws1.Range("J137").Value = 3574.6428
$ws1.Range("K137").Value = 4365.2856
$ws1.Range("L137").Value = 10723.9284
$ws1.Range("M137").Value = -1815.2856
$ws1.Range("N137").Value = -15823.9284

# hunk 5: ARM row 5
$ws2.Range("H5").Value = 153.8
$ws2.Range("I5").Value = 133.85715
$ws2.Range("J5").Value = 200.33333
$ws2.Range("K5").Value = 133.85715
$ws2.Range("L5").Value = 200.33333
$ws2.Range("M5").Value = -21.85714999999999
$ws2.Range("N5").Value = -424.33333

# hunk 6: ARM row 25
$ws2.Range("H25").Value = 1737.5
$ws2.Range("I25").Value = 1000
$ws2.Range("J25").Value = 1983.3334
$ws2.Range("K25").Value = 1000
$ws2.Range("L25").Value = 1983.3334
$ws2.Range("M25").Value = -598
$ws2.Range("N25").Value = -2787.3334

# hunk 7: ARM row 32
$ws2.Range("H32").Value = 4020.3235
$ws2.Range("I32").Value = 3161.6206
$ws2.Range("J32").Value = 9000.799999999999
$ws2.Range("K32").Value = 3161.6206
$ws2.Range("L32").Value = 9000.799999999999
$ws2.Range("M32").Value = -2874.6206
$ws2.Range("N32").Value = -9574.799999999999

# hunk 8: ARM row 110
$ws2.Range("H110").Value = 5543.5625
$ws2.Range("J110").Value = 10270.667
$ws2.Range("L110").Value = 10270.667
$ws2.Range("N110").Value = -14360.667

# hunk 9: ARM row 132
$ws2.Range("H132").Value = 3436.1155
$ws2.Range("I132").Value = 1392.4
$ws2.Range("K132").Value = 4177.200000000001
$ws2.Range("M132").Value = -1647.200000000001

# hunk 10: ARM row 139
$ws2.Range("H139").Value = 88124.625
$ws2.Range("J139").Value = 88124.625
$ws2.Range("L139").Value = 88124.625
$ws2.Range("N139").Value = -98404.625

# hunk 11: BSM row 4
$ws3.Range("H4").Value = 153.8
$ws3.Range("I4").Value = 133.85715
$ws3.Range("J4").Value = 200.33333
$ws3.Range("K4").Value = 133.85715
$ws3.Range("L4").Value = 200.33333
$ws3.Range("M4").Value = -18.85714999999999
$ws3.Range("N4").Value = -430.33333

# hunk 12: BSM row 22
$ws3.Range("H22").Value = 2644.5
$ws3.Range("I22").Value = 290.33334
$ws3.Range("K22").Value = 290.33334
$ws3.Range("M22").Value = -117.33334

# hunk 13: BSM row 81
$ws3.Range("H81").Value = 28858.2
$ws3.Range("J81").Value = 28858.2
$ws3.Range("L81").Value = 28858.2
$ws3.Range("N81").Value = -30980.2

# hunk 14: BSM row 84
$ws3.Range("H84").Value = 28858.2
$ws3.Range("J84").Value = 28858.2
$ws3.Range("L84").Value = 86574.60000000001
$ws3.Range("N84").Value = -97182.60000000001

# hunk 15: BSM row 99
$ws3.Range("H99").Value = 3000
$ws3.Range("I99").Value = 3000
$ws3.Range("K99").Value = 3000
$ws3.Range("M99").Value = -1502

# hunk 16: BSM row 105
$ws3.Range("H105").Value = 21024.5
$ws3.Range("I105").Value = 22945.111
$ws3.Range("K105").Value = 22945.111
$ws3.Range("M105").Value = -21198.111

# hunk 17: BSM row 134
$ws3.Range("H134").Value = 4518.433
$ws3.Range("I134").Value = 2457.3914
$ws3.Range("K134").Value = 7372.174199999999
$ws3.Range("M134").Value = -4837.174199999999

# hunk 18: CRP row 7
$ws4.Range("H7").Value = 599
$ws4.Range("I7").Value = 599
$ws4.Range("K7").Value = 599
$ws4.Range("M7").Value = -486

# hunk 19: CRP row 31
$ws4.Range("H31").Value = 20709.934
$ws4.Range("J31").Value = 28654.596
$ws4.Range("L31").Value = 28654.596
$ws4.Range("N31").Value = -29244.596

# hunk 20: CRP row 34
$ws4.Range("H34").Value = 20709.934
$ws4.Range("J34").Value = 28654.596
$ws4.Range("L34").Value = 28654.596
$ws4.Range("N34").Value = -29058.596

# hunk 21: CRP row 39
$ws4.Range("H39").Value = 14715.25
$ws4.Range("I39").Value = 14620.333
$ws4.Range("J39").Value = 15000
$ws4.Range("K39").Value = 14620.333
$ws4.Range("L39").Value = 15000
$ws4.Range("M39").Value = -14229.333
$ws4.Range("N39").Value = -15782

# hunk 22: CRP row 49
$ws4.Range("H49").Value = 14715.25
$ws4.Range("I49").Value = 14620.333
$ws4.Range("J49").Value = 15000
$ws4.Range("K49").Value = 14620.333
$ws4.Range("L49").Value = 15000
$ws4.Range("M49").Value = -14438.333
$ws4.Range("N49").Value = -15364

# hunk 23: CRP row 94
$ws4.Range("H94").Value = 1607.3
$ws4.Range("I94").Value = 1752.4286
$ws4.Range("J94").Value = 1529.1538
$ws4.Range("K94").Value = 1752.4286
$ws4.Range("L94").Value = 1529.1538
$ws4.Range("M94").Value = -1301.4286
$ws4.Range("N94").Value = -2431.1538

# hunk 24: CUL row 50
$ws5.Range("H50").Value = 20839790
$ws5.Range("I50").Value = 2405
$ws5.Range("J50").Value = 35723636
$ws5.Range("K50").Value = 7215
$ws5.Range("L50").Value = 107170908
$ws5.Range("M50").Value = -6734
$ws5.Range("N50").Value = -107171870

# hunk 25: CUL row 53
$ws5.Range("H53").Value = 20839790
$ws5.Range("I53").Value = 2405
$ws5.Range("J53").Value = 35723636
$ws5.Range("K53").Value = 7215
$ws5.Range("L53").Value = 107170908
$ws5.Range("M53").Value = -6734
$ws5.Range("N53").Value = -107171870

# hunk 26: CUL row 68
$ws5.Range("H68").Value = 3576
$ws5.Range("J68").Value = 3670.3635
$ws5.Range("L68").Value = 11011.0905
$ws5.Range("N68").Value = -12633.0905

# hunk 27: CUL row 71
$ws5.Range("H71").Value = 3576
$ws5.Range("J71").Value = 3670.3635
$ws5.Range("L71").Value = 33033.2715
$ws5.Range("N71").Value = -41145.2715

# hunk 28: CUL row 113
$ws5.Range("H113").Value = 915.0625
$ws5.Range("I113").Value = 576.125
$ws5.Range("J113").Value = 1254
$ws5.Range("K113").Value = 1728.375
$ws5.Range("L113").Value = 3762
$ws5.Range("M113").Value = 441.625
$ws5.Range("N113").Value = -8102

# hunk 29: CUL row 134
$ws5.Range("H134").Value = 7880
$ws5.Range("I134").Value = 7880
$ws5.Range("K134").Value = 23640
$ws5.Range("M134").Value = -18570

# hunk 30: CUL row 138
$ws5.Range("H138").Value = 5591.923
$ws5.Range("J138").Value = 9192.333000000001
$ws5.Range("L138").Value = 27576.999
$ws5.Range("N138").Value = -37856.999

# hunk 31: CUL row 140
$ws5.Range("H140").Value = 911.4666999999999
$ws5.Range("I140").Value = 911.4666999999999
$ws5.Range("K140").Value = 2734.4001
$ws5.Range("M140").Value = 2445.5999

# hunk 32: GSM row 97
$ws6.Range("H97").Value = 1619.8889
$ws6.Range("I97").Value = 1384.5834
$ws6.Range("K97").Value = 1384.5834
$ws6.Range("M97").Value = -888.5834

# hunk 33: LTW row 42
$ws7.Range("H42").Value = 40025
$ws7.Range("I42").Value = 40025
$ws7.Range("J42").Value = 0
$ws7.Range("K42").Value = 40025
$ws7.Range("L42").Value = 0
$ws7.Range("M42").ClearContents()
$ws7.Range("N42").Value = -39462

# hunk 34: LTW row 46
$ws7.Range("H46").Value = 3853.923
$ws7.Range("I46").Value = 1159.8
$ws7.Range("K46").Value = 1159.8
$ws7.Range("M46").Value = -971.8

# hunk 35: LTW row 49
$ws7.Range("H49").Value = 40025
$ws7.Range("I49").Value = 40025
$ws7.Range("J49").Value = 0
$ws7.Range("K49").Value = 40025
$ws7.Range("L49").Value = 0
$ws7.Range("M49").ClearContents()
$ws7.Range("N49").Value = -39878

# hunk 36: LTW row 55
$ws7.Range("H55").Value = 2381604.2
$ws7.Range("I55").Value = 2778435.5
$ws7.Range("J55").Value = 616.3333
$ws7.Range("K55").Value = 2778435.5
$ws7.Range("L55").Value = 616.3333
$ws7.Range("M55").Value = -2778262.5
$ws7.Range("N55").Value = -962.3333

# hunk 37: LTW row 61
$ws7.Range("H61").Value = 3514.125
$ws7.Range("I61").Value = 3157.4443
$ws7.Range("K61").Value = 3157.4443
$ws7.Range("M61").Value = -2955.4443

# hunk 38: LTW row 93
$ws7.Range("H93").Value = 41253.816
$ws7.Range("I93").Value = 48699.75
$ws7.Range("K93").Value = 48699.75
$ws7.Range("M93").Value = -47451.75

# hunk 39: LTW row 113
$ws7.Range("H113").Value = 3514.125
$ws7.Range("I113").Value = 3157.4443
$ws7.Range("K113").Value = 3157.4443
$ws7.Range("M113").Value = -987.4443000000001

# hunk 40: LTW row 132
$ws7.Range("H132").Value = 5346.778
$ws7.Range("I132").Value = 4419
$ws7.Range("J132").Value = 6345.923
$ws7.Range("K132").Value = 13257
$ws7.Range("L132").Value = 19037.769
$ws7.Range("M132").Value = -10727
$ws7.Range("N132").Value = -24097.769

# hunk 41: LTW row 136
$ws7.Range("H136").Value = 8068.3257
$ws7.Range("I136").Value = 6486.4375
$ws7.Range("K136").Value = 19459.3125
$ws7.Range("M136").Value = -16909.3125

# hunk 42: WVR row 122
$ws8.Range("H122").Value = 5518
$ws8.Range("I122").Value = 2374.6155
$ws8.Range("K122").Value = 7123.8465
$ws8.Range("M122").Value = -4673.8465

# hunk 43: WVR row 126
$ws8.Range("H126").Value = 4157.933
$ws8.Range("I126").Value = 2930.9473
$ws8.Range("K126").Value = 8792.841899999999
$ws8.Range("M126").Value = -6322.841899999999

# hunk 44: WVR row 132
$ws8.Range("H132").Value = 6974
$ws8.Range("I132").Value = 5218.1
$ws8.Range("J132").Value = 9900.5
$ws8.Range("K132").Value = 15654.3
$ws8.Range("L132").Value = 29701.5
$ws8.Range("M132").Value = -13124.3
$ws8.Range("N132").Value = -34761.5
